$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "Handoff failed" status text (shared string) -> "Not yet
#    handed off". This string is referenced from the Overview sheet (B2/C2)
#    and from the per-language sheets' Status column (B2), so editing the
#    cells that hold it updates every occurrence uniformly.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Not yet handed off"
$overview.Range("C2").Value = "Not yet handed off"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = "Not yet handed off"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = "Not yet handed off"

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: a handoff report was generated for the .md file.
#    - a new hyperlink to the generated .xlf handoff package goes in C2
#      (A2's and A3's existing hyperlinks are left exactly as they were)
#    - the handoff datetime (D2) is recorded
#    - the handoff reason (H2) switches from "Ignored" to "Include"
# ---------------------------------------------------------------------------
$zhcnXlf = "22cb96f6-bddf-4d90-b60e-c0e52b6f4bf7.6dec47758dcd162a0333e842f0647c0269c4c04f.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/96aae2cc95ed185ab62e33bd34915a18d2752088/e2e/.loc/$zhcnXlf", "", "", $zhcnXlf)

$zhcn.Range("D2").Value = "2016-01-08 14:53:28"
$zhcn.Range("H2").Value = "Include"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same handoff report, for the de-de target language.
# ---------------------------------------------------------------------------
$dedeXlf = "22cb96f6-bddf-4d90-b60e-c0e52b6f4bf7.6dec47758dcd162a0333e842f0647c0269c4c04f.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/96aae2cc95ed185ab62e33bd34915a18d2752088/e2e/.loc/$dedeXlf", "", "", $dedeXlf)

$dede.Range("D2").Value = "2016-01-08 14:53:42"
$dede.Range("H2").Value = "Include"
